# Mark several "Programming Status" cells (column D) as Done,
# matching the green "Done" formatting already used elsewhere in the sheet (e.g. C4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$doneStyleSource = $ws.Range("C4")
$doneColor = $doneStyleSource.Interior.Color

$targetCells = @("D4", "D5", "D7", "D8", "D9", "D10", "D11")

foreach ($addr in $targetCells) {
    $cell = $ws.Range($addr)
    $cell.Value = "Done"
    $cell.Interior.Color = $doneColor
    $cell.Font.Bold = $doneStyleSource.Font.Bold
    $cell.Font.Name = $doneStyleSource.Font.Name
    $cell.Font.Size = $doneStyleSource.Font.Size
}
